# Applies "updated raw data with spawn/brood data" commit to the
# "data_table_factors - data" worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data_table_factors - data")

# ---------------------------------------------------------------------------
# 1. Updated paper / author / species text (shared strings 164/165/166),
#    plus the associated new rows 328-339 that now carry A (paper) and
#    B (author) values, and C328:C333 species correction.
# ---------------------------------------------------------------------------
$paperTitle  = "Large birth size does not reduce negative latent effects of harsh environments across life stages in two coral species"
$authorList  = "Hartman, Marhaver, Chamberland, Sandin and Vermeij"
$speciesName = "Agaricia humlis"

$ws.Range("A328:A339").Value = $paperTitle
$ws.Range("B328:B339").Value = $authorList
$ws.Range("C328:C333").Value = $speciesName

# ---------------------------------------------------------------------------
# 2. New "spawn/brood" column E values.
#    E319 previously held the (now repurposed) paper-title string with a
#    bold-ish direct style (s="2") -- clear that formatting before writing
#    the new "spawn" value so the stray cellXf stops being referenced.
# ---------------------------------------------------------------------------
$ws.Range("E319").ClearFormats()

$ws.Range("E69:E73").Value = "spawn"
$ws.Range("E79:E83").Value = "spawn"
$ws.Range("E106:E110").Value = "spawn"
$ws.Range("E184:E327").Value = "spawn"
$ws.Range("E334:E339").Value = "spawn"

$ws.Range("E328:E333").Value = "brood"

# ---------------------------------------------------------------------------
# 3. D334:D339 factor/treatment index bumped from 57 to 58.
# ---------------------------------------------------------------------------
$ws.Range("D334:D339").Value = 58

# ---------------------------------------------------------------------------
# 4. Column C (species) width now has an explicit, best-fit width.
# ---------------------------------------------------------------------------
$ws.Columns(3).ColumnWidth = 19

# ---------------------------------------------------------------------------
# 5. View state: unfreeze/refreeze at the top row (instead of row 322) and
#    move the active selection to E333.
# ---------------------------------------------------------------------------
$ws.Activate()
$excel.ActiveWindow.FreezePanes = $false
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("E333").Select()
